$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.259.87'
$ws.Range('E2').Value = '  +0.72%  '
$ws.Range('D3').Value = '2.265.14'
$ws.Range('E3').Value = '  +1.27%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '497.15'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.82%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '128.93'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.71%  '
$ws.Range('E7').Value = '  +0.54%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.525'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.0953'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.67%  '
$ws.Range('E10').Value = '  +1.30%  '
$ws.Range('E11').Value = '  +4.20%  '
$ws.Range('E12').Value = '  +3.41%  '
$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '22.92'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +5.56%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '2.666.50'
$ws.Range('E14').Value = '  +0.95%  '
$ws.Range('D15').Value = '54.249.01'
$ws.Range('E15').Value = '  +0.75%  '
$ws.Range('E16').Value = '  +0.85%  '
$ws.Range('D17').Value = '2.270.04'
$ws.Range('E17').Value = '  +0.91%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '10.24'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +2.19%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.15'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.71%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '303.50'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.98%  '
$ws.Range('E21').Value = '  -2.17%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.999'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.46%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '61.00'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.95%  '
$ws.Range('E24').Value = '  -0.79%  '
$ws.Range('E25').Value = '  +0.41%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.34'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +3.99%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '171.03'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.58%  '
$ws.Range('E28').Value = '  +1.04%  '
$ws.Range('D29').Value = '0.0₃0693'
$ws.Range('E29').Value = '  +1.38%  '
$ws.Range('E30').Value = '  +1.15%  '
$ws.Range('E31').Value = '  +1.46%  '
$ws.Range('E32').Value = '  +0.14%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '17.80'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.33%  '
$ws.Range('E34').Value = '  +0.55%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.936'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +8.97%  '
$ws.Range('E36').Value = '  +0.59%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.71'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.59%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.375'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.59%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.39'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.31%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.37'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.97%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '125.08'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -2.21%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '4.81'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.82%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0493'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +2.34%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0892'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.73%  '
$ws.Range('E45').Value = '  +1.16%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '241.85'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.90%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.374'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.55%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0205'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.76%  '
$ws.Range('E49').Value = '  +0.49%  '
$ws.Range('E50').Value = '  -0.02%  '
$ws.Range('E51').Value = '  -0.55%  '
